$wb = $excel.ActiveWorkbook


$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1144751
$ws.Range("I19").Value = 2024735.4
$ws.Range("J19").Value = 771.2
$ws.Range("K19").Value = 2024735.4
$ws.Range("L19").Value = 771.2
$ws.Range("M19").Value = -2024560.4
$ws.Range("N19").Value = -1121.2
$ws.Range("H33").Value = 624.25
$ws.Range("I33").Value = 665.6667
$ws.Range("K33").Value = 665.6667
$ws.Range("M33").Value = -436.6667
$ws.Range("H44").Value = 24375
$ws.Range("J44").Value = 24375
$ws.Range("L44").Value = 24375
$ws.Range("N44").Value = -25299
$ws.Range("H98").Value = 5317.6665
$ws.Range("I98").Value = 1987.4286
$ws.Range("J98").Value = 9980
$ws.Range("K98").Value = 1987.4286
$ws.Range("L98").Value = 9980
$ws.Range("M98").Value = -489.4286
$ws.Range("N98").Value = -12976
$ws.Range("H99").Value = 718.1667
$ws.Range("I99").Value = 785.6
$ws.Range("J99").Value = 381
$ws.Range("K99").Value = 2356.8
$ws.Range("L99").Value = 1143
$ws.Range("M99").Value = -858.8000000000002
$ws.Range("N99").Value = -4139
$ws.Range("H112").Value = 470284.84
$ws.Range("I112").Value = 817.5
$ws.Range("J112").Value = 521038.06
$ws.Range("K112").Value = 2452.5
$ws.Range("L112").Value = 1563114.18
$ws.Range("M112").Value = -1344.5
$ws.Range("N112").Value = -1565330.18
$ws.Range("H118").Value = 898.61536
$ws.Range("I118").Value = 762.1667
$ws.Range("J118").Value = 1015.5714
$ws.Range("K118").Value = 2286.5001
$ws.Range("L118").Value = 3046.7142
$ws.Range("M118").Value = -629.5001000000002
$ws.Range("N118").Value = -6360.7142
$ws.Range("H122").Value = 5317.6665
$ws.Range("I122").Value = 1987.4286
$ws.Range("J122").Value = 9980
$ws.Range("K122").Value = 5962.2858
$ws.Range("L122").Value = 29940
$ws.Range("M122").Value = -3512.2858
$ws.Range("N122").Value = -34840
$ws.Range("H137").Value = 4037.3333
$ws.Range("I137").Value = 3418.3333
$ws.Range("J137").Value = 4346.8335
$ws.Range("K137").Value = 10254.9999
$ws.Range("L137").Value = 13040.5005
$ws.Range("M137").Value = -7704.999899999999
$ws.Range("N137").Value = -18140.5005
$ws.Range("H138").Value = 2391.74
$ws.Range("I138").Value = 1104.4615
$ws.Range("J138").Value = 2584.092
$ws.Range("K138").Value = 3313.3845
$ws.Range("L138").Value = 7752.276
$ws.Range("M138").Value = 1826.6155
$ws.Range("N138").Value = -18032.276
$ws.Range("H141").Value = 7925.7095
$ws.Range("I141").Value = 8274.893
$ws.Range("J141").Value = 4666.6665
$ws.Range("K141").Value = 24824.679
$ws.Range("L141").Value = 13999.9995
$ws.Range("M141").Value = -19644.679
$ws.Range("N141").Value = -24359.9995

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1456
$ws.Range("I45").Value = 1440.3334
$ws.Range("J45").Value = 1471.6666
$ws.Range("K45").Value = 1440.3334
$ws.Range("L45").Value = 1471.6666
$ws.Range("M45").Value = -1063.3334
$ws.Range("N45").Value = -2225.6666
$ws.Range("H61").Value = 957.12
$ws.Range("I61").Value = 957.12
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 957.12
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -745.12
$ws.Range("N61").Value = ""
$ws.Range("H74").Value = 5315.95
$ws.Range("I74").Value = 5012.1665
$ws.Range("J74").Value = 8050
$ws.Range("K74").Value = 5012.1665
$ws.Range("L74").Value = 8050
$ws.Range("M74").Value = -4138.1665
$ws.Range("N74").Value = -9798
$ws.Range("H77").Value = 5315.95
$ws.Range("I77").Value = 5012.1665
$ws.Range("J77").Value = 8050
$ws.Range("K77").Value = 25060.8325
$ws.Range("L77").Value = 40250
$ws.Range("M77").Value = -20692.8325
$ws.Range("N77").Value = -48986
$ws.Range("H102").Value = 1994.7858
$ws.Range("I102").Value = 1912.7
$ws.Range("K102").Value = 1912.7
$ws.Range("M102").Value = -290.7
$ws.Range("H122").Value = 2950.3333
$ws.Range("I122").Value = 1820
$ws.Range("J122").Value = 4532.8
$ws.Range("K122").Value = 5460
$ws.Range("L122").Value = 13598.4
$ws.Range("M122").Value = -3010
$ws.Range("N122").Value = -18498.4
$ws.Range("H132").Value = 1651.0212
$ws.Range("I132").Value = 994.87805
$ws.Range("J132").Value = 6134.6665
$ws.Range("K132").Value = 2984.63415
$ws.Range("L132").Value = 18403.9995
$ws.Range("M132").Value = -454.6341499999999
$ws.Range("N132").Value = -23463.9995
$ws.Range("H136").Value = 957.12
$ws.Range("I136").Value = 957.12
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 2871.36
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -321.3600000000001
$ws.Range("N136").Value = ""

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").Value = ""
$ws.Range("H94").Value = 1601.4286
$ws.Range("I94").Value = 1975
$ws.Range("K94").Value = 1975
$ws.Range("M94").Value = -1524
$ws.Range("H134").Value = 1663.909
$ws.Range("I134").Value = 1017.7353
$ws.Range("J134").Value = 3860.9
$ws.Range("K134").Value = 3053.2059
$ws.Range("L134").Value = 11582.7
$ws.Range("M134").Value = -518.2058999999999
$ws.Range("N134").Value = -16652.7

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 71434520
$ws.Range("I31").Value = 831
$ws.Range("J31").Value = 100008000
$ws.Range("K31").Value = 831
$ws.Range("L31").Value = 100008000
$ws.Range("M31").Value = -536
$ws.Range("N31").Value = -100008590
$ws.Range("H34").Value = 71434520
$ws.Range("I34").Value = 831
$ws.Range("J34").Value = 100008000
$ws.Range("K34").Value = 831
$ws.Range("L34").Value = 100008000
$ws.Range("M34").Value = -629
$ws.Range("N34").Value = -100008404
$ws.Range("H58").Value = 1751.7089
$ws.Range("I58").Value = 1583.409
$ws.Range("J58").Value = 2606.1538
$ws.Range("K58").Value = 1583.409
$ws.Range("L58").Value = 2606.1538
$ws.Range("M58").Value = -1380.409
$ws.Range("N58").Value = -3012.1538
$ws.Range("H132").Value = 2900.1592
$ws.Range("I132").Value = 2309.484
$ws.Range("J132").Value = 4308.6924
$ws.Range("K132").Value = 6928.451999999999
$ws.Range("L132").Value = 12926.0772
$ws.Range("M132").Value = -4398.451999999999
$ws.Range("N132").Value = -17986.0772
$ws.Range("H134").Value = 4311
$ws.Range("I134").Value = 4892.4165
$ws.Range("J134").Value = 2760.5557
$ws.Range("K134").Value = 14677.2495
$ws.Range("L134").Value = 8281.667099999999
$ws.Range("M134").Value = -12142.2495
$ws.Range("N134").Value = -13351.6671
$ws.Range("H136").Value = 1751.7089
$ws.Range("I136").Value = 1583.409
$ws.Range("J136").Value = 2606.1538
$ws.Range("K136").Value = 4750.227000000001
$ws.Range("L136").Value = 7818.4614
$ws.Range("M136").Value = -2200.227000000001
$ws.Range("N136").Value = -12918.4614

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1375
$ws.Range("J68").Value = 1166.6666
$ws.Range("L68").Value = 3499.9998
$ws.Range("N68").Value = -5121.9998
$ws.Range("H71").Value = 1375
$ws.Range("J71").Value = 1166.6666
$ws.Range("L71").Value = 10499.9994
$ws.Range("N71").Value = -18611.9994
$ws.Range("H113").Value = 531.6316
$ws.Range("I113").Value = 527.3095
$ws.Range("J113").Value = 543.73334
$ws.Range("K113").Value = 1581.9285
$ws.Range("L113").Value = 1631.20002
$ws.Range("M113").Value = 588.0715
$ws.Range("N113").Value = -5971.20002
$ws.Range("H116").Value = 3746.6667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1918.4375
$ws.Range("I102").Value = 1488.9286
$ws.Range("K102").Value = 1488.9286
$ws.Range("M102").Value = 133.0714
$ws.Range("H122").Value = 6737.5
$ws.Range("I122").Value = 3983.3333
$ws.Range("J122").Value = 15000
$ws.Range("K122").Value = 11949.9999
$ws.Range("L122").Value = 45000
$ws.Range("M122").Value = -9499.999899999999
$ws.Range("N122").Value = -49900
$ws.Range("H132").Value = 1532.4565
$ws.Range("I132").Value = 983.2143
$ws.Range("J132").Value = 7299.5
$ws.Range("K132").Value = 2949.6429
$ws.Range("L132").Value = 21898.5
$ws.Range("M132").Value = -419.6428999999998
$ws.Range("N132").Value = -26958.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 6010
$ws.Range("I122").Value = 2770
$ws.Range("K122").Value = 8310
$ws.Range("M122").Value = -5860
$ws.Range("H132").Value = 8652.9375
$ws.Range("I132").Value = 10216.357
$ws.Range("J132").Value = 6464.15
$ws.Range("K132").Value = 30649.071
$ws.Range("L132").Value = 19392.45
$ws.Range("M132").Value = -28119.071
$ws.Range("N132").Value = -24452.45
$ws.Range("H136").Value = 3130.0833
$ws.Range("I136").Value = 782.625
$ws.Range("J136").Value = 7825
$ws.Range("K136").Value = 2347.875
$ws.Range("L136").Value = 23475
$ws.Range("M136").Value = 202.125
$ws.Range("N136").Value = -28575

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6411740.5
$ws.Range("I132").Value = 872.6316
$ws.Range("J132").Value = 23812666
$ws.Range("K132").Value = 2617.8948
$ws.Range("L132").Value = 71437998
$ws.Range("M132").Value = -87.89480000000003
$ws.Range("N132").Value = -71443058
$ws.Range("H136").Value = 1663.1034
$ws.Range("I136").Value = 1178.6364
$ws.Range("K136").Value = 3535.9092
$ws.Range("M136").Value = -985.9092000000001
